$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.535.49'
$ws.Range("E2").Value = '  -2.20%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.993.95'
$ws.Range("E3").Value = '  -0.74%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.88'
$ws.Range("E5").Value = '  -9.89%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.600'
$ws.Range("E6").Value = '  -2.44%  '

$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.96'
$ws.Range("E8").Value = '  -1.70%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.370'
$ws.Range("E9").Value = '  -3.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.02'
$ws.Range("E10").Value = '  +2.78%  '

$ws.Range("E11").Value = '  -2.61%  '

$ws.Range("E12").Value = '  -3.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.19'
$ws.Range("E13").Value = '  -0.13%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.290.01'
$ws.Range("E14").Value = '  -0.60%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.28'
$ws.Range("E15").Value = '  -3.20%  '

$ws.Range("E16").Value = '  -5.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.08'
$ws.Range("E17").Value = '  -2.76%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.002.60'
$ws.Range("E18").Value = '  -2.43%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '36.491.75'
$ws.Range("E19").Value = '  -2.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '67.72'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0803'
$ws.Range("E21").Value = '  -3.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.29'
$ws.Range("E22").Value = '  +3.45%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '221.90'
$ws.Range("E23").Value = '  -2.55%  '

$ws.Range("E24").Value = '  -0.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.36'
$ws.Range("E25").Value = '  +0.61%  '

$ws.Range("E26").Value = '  -8.15%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.93'
$ws.Range("E27").Value = '  -1.67%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.66'
$ws.Range("E28").Value = '  -2.11%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.126'
$ws.Range("E29").Value = '  -3.11%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '18.80'
$ws.Range("E30").Value = '  -3.91%  '

$ws.Range("E31").Value = '  +0.88%  '

$ws.Range("E32").Value = '  -2.45%  '

$ws.Range("E33").Value = '  -5.02%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0602'
$ws.Range("E34").Value = '  -6.16%  '

$ws.Range("E35").Value = '  -5.73%  '

$ws.Range("E36").Value = '  -1.14%  '

$ws.Range("E37").Value = '  +0.18%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.35'
$ws.Range("E38").Value = '  +0.82%  '

$ws.Range("E39").Value = '  -2.81%  '

$ws.Range("E40").Value = '  +6.73%  '

$ws.Range("E41").Value = '  -1.39%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0939'
$ws.Range("E42").Value = '  +1.13%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.457.08'
$ws.Range("E43").Value = '  +4.13%  '

$ws.Range("E44").Value = '  -3.84%  '

$ws.Range("E45").Value = '  -8.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.03'
$ws.Range("E46").Value = '  -0.33%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.17'
$ws.Range("E47").Value = '  -3.34%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.990'
$ws.Range("E48").Value = '  -2.84%  '

$ws.Range("E49").Value = '  -0.78%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.82'
$ws.Range("E50").Value = '  -2.91%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.70'
$ws.Range("E51").Value = '  +7.44%  '
